$d = $word.ActiveDocument

# --- Paragraph 1: title block ---
# Before: <w:br/><w:t>date (22.01.25)</w:t><w:br/><w:t>MONOFORMER title</w:t>
# After:  <w:t>date (21.01.25)</w:t><w:br/><w:t>Time-MoE title</w:t>  (leading break removed)
$p1 = $d.Paragraphs(1).Range
$expectedOld = [char]11 + "המאמר היומי של מייק - 22.01.25" + [char]11 + "MONOFORMER: ONE TRANSFORMER FOR BOTH DIFFUSION AND AUTOREGRESSION" + [char]13
if ($p1.Text -ne $expectedOld) { throw "Paragraph 1 did not match expected original content" }
$p1.Text = "המאמר היומי של מייק - 21.01.25" + [char]11 + "Time-MoE: Billion-Scale Time Series Foundation Models with Mixture of Experts"

# --- Paragraph 2: Intro paragraph -> now about Time-MoE / time-series foundation models ---
$rng = $d.Content
$found = $rng.Find.Execute("היום נעשה סקירה קצרה של מאמר די מעניין ששילב שני סוגים של מודלים, מודל שפה ומודל ויז'ן בטרנספורמר אחד. רוב המודלים מולטימודליים מורכבים מכמה מודלים שכל אחד מהם אחראי על הגנרוט של סוג דאטה אחד. למשל מודלי שפה ויזואליים בד״כ מורכבים משני מודלים: מודל שפה ומודל לגנרוט תמונות. המחברים מציעים ״לחבר״ את שני המודלים האלה למודל טרנספורמר אחד וזה נעשה בצורה די אינטואיטיבית.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Paragraph 2 source text not found" }
$rng.Text = "המאמר משך את תשומת ליבי למרות הידע הרדוד שאני מחזיק לגבי תחום הסדרות העתיות (time-series). בגדול הסיבה העיקרית לכך שבשמו מופיע צמד מילים ""Foundational Models"" שזה חיה די נדירה בתחום הסדרות העתיות להבדיל מתחום מודלי שפה. הסיבה לכך (כנראה) היא מגוון עשיר הרבה יותר של סדרות עתיות השונות יחסית לשפה טבעית."

# --- Paragraph 3: Architecture paragraph -> SwiGLU token embedding discussion ---
$rng = $d.Content
$found = $rng.Find.Execute("קודם כל נציין כי שני המודלים האלו עובדים במרחב הטוקנים כאשר עבור מודלי שפה כל טוקן הוא חלק של מילה או מילה שלמה ואילו עבור מודל ויזואלי כל טוקן הוא פאץ' של תמונה. אז הניסיון לחבר אותם למודל אחד נראה די טבעי אך לא ברור האם ניתן לאמן אותו הטרנספורמר לגנרט שפה ותמונות כאחד. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Paragraph 3 source text not found" }
$rng.Text = "האמת לא מצאתי ב- Time-MoE, המבוססת כמובן על הטרנספורמרים, מציאות ארכיטקטוניות מאוד מעניינות ועם זאת יש בו כמה דברים שונים מאלו שאנו רגילים לראות ב-LLMs. למשל במקום שכבת טוקניזציה ואמבדינג, מבוססים על מילון טוקנים, שיש לנו ב-LLMs במודל המוצע יש כל טוקן (שזו נקודה בסדרה) עובר טרנספורמציה לא לינאריות עם אקטיבציה מסוג SwiGLU וכמה טרנספורמציות לינאריות."

# --- Paragraph 4: Transformer-layer paragraph -> MoE / RMSNorm discussion ---
$rng = $d.Content
$found = $rng.Find.Execute("המודל המוצע מגנרט שפה בדיוק כמו LLM רגיל, בצורה אוטורגרסיבית, כלומר, טוקן אחרי טוקן. אבל איך ניתן לשלב אותו עם מודל לגנרוט תמונות שכמובן מבוסס על מודלי דיפוזיה (בשנת 2025 זה האופציה הדיפולטית הרי). קודם כל צריך לזכור שמודל אוטורגרסיבי (לגנרוט שפה) עובד בצורה סיבתית (קוזלית), כלומר במהלך גנרוט טוקן n כל הטוקנים מאחוריו ממוסכים ולא משתתפים בגנרוט(משתמשים במסכה קוזלית). למודלי אנו צריכים מודל דו כיווני כי בזמן גנרוט פאץ' של תמונה כדאי מאוד להשתמש בכל הפאצ'ים האחרים. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Paragraph 4 source text not found" }
$rng.Text = "בנוגע לשכבת הטרנספורמרים, המחברים לוקחים ארכיטקטורת MoE די סטנדרטית. השוני היחיד שמשך את עיניי הוא שימוש בשיטת נרמול RMSNorm שלא הכרתי. פרט לכך יש את כל השכבות הרגילות של הטרנספורמרים כולל כמובן שכבות residual."

# --- Paragraph 5: Output-layer paragraph -> multi-head forecasting discussion ---
$rng = $d.Content
$found = $rng.Find.Execute("בדיוק כך בנוי המודל המוצע - השפה מגונרטת עם מסכה קוזלית והתמונה מגונרטת עם כל הטוקנים (כולל הטוקנים של טקסט). דרך אגב הגישה הזו תעבוד גם לכיוון השני: כלומר בגנרוט של טקסט מתמונה (למשל למשימת captioning). אבל איך נדע לעבור ממצב ״קוזלי״ למצב ״דו-כיווני״. המחברים מציעים להשתמש בטוקן מסוים המסמן שממנו מתחיל גנרוט התמונה - הטוקן הזה אמור להיות מג'ונרט למשל למשימה יצירת תמונה מטקסט.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Paragraph 5 source text not found" }
$rng.Text = "השכבה האחרונה של Time-MoE היא קצת שונה ממה שאנו רגילים לראות בטרנספורמרים. מכיוון שלהבדיל ממודלי שפה אנו צריכים מודל בעולם של TS אנו צריכים לחזות במספר נקודות זמן שונה (נגיד שניה, דקה או יום קדימה), המחברים משתמשים בכמה ראשים בשכבה האחרונה. כל ראש אחראי על חיזוי באופק מסוים (כמות דגימות קדימה). באימון משלבים את הלוסים מכל הראשים."

# --- Paragraph 6: Loss-function paragraph -> Huber loss / load-balancing discussion ---
$rng = $d.Content
$found = $rng.Find.Execute("כמה מילים על הטרנספורמר לגנרוט תמונה. המאמר משתמש במודל דיפוזיה לטנטי כאשר המודל מאומן לבנות ייצוג לטנטי של תמונה מרעש (עבור כל פאץ). לאחר מכן כל הייצוגים (של הפאצ'ים) מועברים דרך הדקודר (מבוסס VAE) שבונה ממנו תמונה. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Paragraph 6 source text not found" }
$rng.Text = "גם פונקציות לוס במאמר הן די סטנדרטיות: פונקצית הובר שהיא הגרסה הרובסטית של L2 (הלא נותנת לא להגיע לערכים גבוהים מאוד). בנוסף יש איבר רגולריזציה שמנסה להפעיל את כל המומחים ב-MoE בצורה אחידה. וכמובן אימנו את המודל על דאטהסטים ענקיים ומגוונים."

# --- Paragraph 7: Closing paragraph -> short sign-off ---
$rng = $d.Content
$found = $rng.Find.Execute("המודל מאומן עם הלוס שהוא סכום משוקלל של הלוסים הסטנדרטיים עבור המודלים המוזכרים: מודל שפה ומודל דיפוזיה. המאמר מצליח לגנרט תמונות די יפות….", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Paragraph 7 source text not found" }
$rng.Text = "זהו וזה - סקירה קצרה, ובתקווה גם ברורה…. "

# --- Paragraph 8: arXiv link paragraph -> updated link ---
$rng = $d.Content
$found = $rng.Find.Execute("https://arxiv.org/abs/2409.16280", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Paragraph 8 source text not found" }
$rng.Text = "https://arxiv.org/pdf/2409.16040"

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
Write-Output "All replacements applied successfully."
